$wb = $excel.ActiveWorkbook

# New salesperson to insert right before "SALAZAR VERA ENRIQUE WILLIAM"
$newName = "RAMIREZ MOREIRA MAYRA JACQUELINE"
$office = "OFICINA-CATAECSA"

# --- Sheet "VENTAS POR GRUPO": 18 columns (A..R), summary row uses "X de 44" -> "X de 45"
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Rows.Item(39).Insert()
$ws1.Range("A39").Value = $office
$ws1.Range("B39").Value = $newName
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(39, $col).Value = 0
}
# Fix the trailing summary row text counts ("X de 44" -> "X de 45") now sitting at row 47
for ($col = 3; $col -le 18; $col++) {
    $cell = $ws1.Cells.Item(47, $col)
    $cell.Value = $cell.Value() -replace "de 44", "de 45"
}

# --- Sheet "VENTA MENSUAL": 7 columns (A..G), summary row keeps numeric totals
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Rows.Item(39).Insert()
$ws2.Range("A39").Value = $office
$ws2.Range("B39").Value = $newName
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(39, $col).Value = 0
}
